# Applies the Coinranking price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions": updates columns D (Price)
# and E (Volume(1h)) for rows 2-51. Values are kept as plain text (matching
# the workbook's original inline-string cells) by forcing a text number
# format before assignment and clearing it again afterwards so no stray
# cell style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @("30.515.39", "  +0.73%  ")
    3 = @("1.871.17", "  +0.19%  ")
    4 = @("0.9993", "  -0.07%  ")
    5 = @("247.07", "  +1.51%  ")
    6 = @("0.9995", "  -0.07%  ")
    7 = @("0.4734", "  +0.18%  ")
    8 = @("0.2917", "  +1.62%  ")
    9 = @("0.06489", "  +0.18%  ")
    10 = @("22.19", "  +5.72%  ")
    11 = @("0.07713", "  -0.23%  ")
    12 = @("97.43", "  +2.35%  ")
    13 = @("0.7435", "  +5.73%  ")
    14 = @("1.868.64", "  -0.08%  ")
    15 = @("5.154", "  +1.07%  ")
    16 = @("273.70", "  -0.53%  ")
    17 = @("30.495.17", "  +0.72%  ")
    18 = @("13.41", "  +0.54%  ")
    19 = @($null, "  +0.03%  ")
    20 = @("0.000007513", "  -0.55%  ")
    21 = @("2.114.54", "  -0.05%  ")
    22 = @("0.9991", "  -0.09%  ")
    23 = @("5.250", "  +0.89%  ")
    24 = @("6.166", "  +0.54%  ")
    25 = @("9.284", "  -0.26%  ")
    26 = @("163.64", "  -1.08%  ")
    27 = @("18.75", "  -1.19%  ")
    28 = @("1.918", "  +0.53%  ")
    29 = @("0.09998", "  +1.43%  ")
    30 = @("1.349", "  -2.06%  ")
    31 = @("1.509", "  -0.12%  ")
    32 = @("4.284", "  +0.86%  ")
    33 = @("4.101", "  +1.80%  ")
    34 = @("0.04810", "  +1.62%  ")
    35 = @("1.118", "  -0.27%  ")
    36 = @("0.6949", "  +0.53%  ")
    37 = @("2.711", "  +0.21%  ")
    38 = @("0.01851", "  +0.48%  ")
    39 = @("2.741", "  +0.07%  ")
    40 = @("6.216", "  -1.93%  ")
    41 = @("72.90", "  +4.01%  ")
    42 = @("1.968", "  +3.70%  ")
    43 = @("0.4187", "  +2.48%  ")
    44 = @("0.9993", "  -0.06%  ")
    45 = @("0.8335", "  -1.03%  ")
    46 = @("102.32", "  +0.29%  ")
    47 = @("9.311", "  +0.75%  ")
    48 = @("35.41", "  +1.39%  ")
    49 = @("6.970", "  -1.34%  ")
    50 = @("924.11", "  +0.23%  ")
    51 = @("0.05638", "  +1.21%  ")
}

foreach ($row in $updates.Keys) {
    $priceValue = $updates[$row][0]
    $volumeValue = $updates[$row][1]

    if ($priceValue -ne $null) {
        $priceCell = $ws.Range("D$row")
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $priceValue
        $priceCell.ClearFormats()
    }

    $ws.Range("E$row").Value = $volumeValue
}
